$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record was added to the "Zapallo" sheet. Insert a fresh row at 118,
# pushing the existing rows 118 (Camote / Region del Maule) and 119 (Paine /
# Region de O'Higgins) down to 119 and 120 respectively, keeping them intact.
$ws.Rows.Item(118).Insert()

# Seed the new row 118 with a copy of the (now shifted) row 119 so every
# column starts out populated and correctly typed/styled, then overwrite the
# handful of cells that actually differ for the new weekly record.
$ws.Range("A119:R119").Copy()
$ws.Range("A118").PasteSpecial()

$ws.Range("D118").Value = 44568
$ws.Range("I118").Value = "1a nueva(o)"
$ws.Range("J118").Value = 400
$ws.Range("K118").Value = 350
$ws.Range("L118").Value = 400
$ws.Range("M118").Value = 375
$ws.Range("O118").Value = "Región de O'Higgins"
$ws.Range("P118").Value = 375
